# Auto-generated cell updates applying the diff to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.873.22"
$ws.Range("E2").Value = "  +4.37%  "
$ws.Range("D3").Value = "3.633.86"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'631.98"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Value = "'160.10"
$ws.Range("E6").Value = "  +5.26%  "
$ws.Range("D7").Value = "3.632.16"
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.496"
$ws.Range("E10").Value = "  +6.72%  "
$ws.Range("D11").Value = "'7.26"
$ws.Range("E11").Value = "  +6.08%  "
$ws.Range("D12").Value = "'0.442"
$ws.Range("E12").Value = "  +3.69%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("E13").Value = "  +5.07%  "
$ws.Range("D14").Value = "'33.47"
$ws.Range("E14").Value = "  +5.80%  "
$ws.Range("D15").Value = "4.246.24"
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "3.628.28"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "69.695.55"
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("E19").Value = "  +5.81%  "
$ws.Range("D20").Value = "'16.02"
$ws.Range("E20").Value = "  +4.28%  "
$ws.Range("D21").Value = "'10.24"
$ws.Range("E21").Value = "  +11.40%  "
$ws.Range("D22").Value = "'464.92"
$ws.Range("E22").Value = "  +4.83%  "
$ws.Range("D23").Value = "'0.645"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("D24").Value = "'78.73"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("E25").Value = "  +12.13%  "
$ws.Range("D26").Value = "'10.79"
$ws.Range("E26").Value = "  +5.68%  "
$ws.Range("D27").Value = "3.774.92"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'9.26"
$ws.Range("E29").Value = "  +13.28%  "
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("D31").Value = "'1.73"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("D32").Value = "'0.179"
$ws.Range("D33").Value = "'6.64"
$ws.Range("E33").Value = "  +8.10%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +6.00%  "
$ws.Range("D36").Value = "'26.61"
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("D37").Value = "3.626.24"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("D38").Value = "'8.48"
$ws.Range("E38").Value = "  +5.96%  "
$ws.Range("D39").Value = "'2.46"
$ws.Range("E39").Value = "  +15.22%  "
$ws.Range("D41").Value = "'0.0936"
$ws.Range("E41").Value = "  +8.93%  "
$ws.Range("D43").Value = "'177.19"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "'5.65"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").Value = "'32.00"
$ws.Range("E45").Value = "  +17.71%  "
$ws.Range("D46").Value = "'0.916"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.84"
$ws.Range("E47").Value = "  +11.36%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.37"
$ws.Range("E48").Value = "  +12.37%  "
$ws.Range("D49").Value = "'46.42"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("E51").Value = "  +8.07%  "
